# Rename the worksheet to the batter's name
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Saurabh Tiwary"

# Insert a new column before column A ("matchNo"), shifting
# teamName..result from A:L to B:M
$ws.Columns.Item(1).Insert()

# New header for column A
$ws.Cells.Item(1, 1).Value = "matchNo"

# matchNo value for the pre-existing data row (now row 2)
$ws.Cells.Item(2, 1).Value = "42nd"

# The numeric-looking stat columns (runs, balls, fours, sixes, sr) must be
# stored as text, like the rest of the sheet - force text format on those
# cells before assigning so Excel doesn't coerce them into numbers.
$statCols = @(5, 6, 7, 8, 9)
foreach ($col in $statCols) {
    $ws.Cells.Item(3, $col).NumberFormat = "@"
    $ws.Cells.Item(4, $col).NumberFormat = "@"
    $ws.Cells.Item(5, $col).NumberFormat = "@"
}

# Row 3: 46th, vs Delhi Capitals
$ws.Cells.Item(3, 1).Value = "46th"
$ws.Cells.Item(3, 2).Value = "Mumbai Indians"
$ws.Cells.Item(3, 3).Value = "Saurabh Tiwary"
$ws.Cells.Item(3, 4).Value = "c †Pant b Patel"
$ws.Cells.Item(3, 5).Value = "15"
$ws.Cells.Item(3, 6).Value = "18"
$ws.Cells.Item(3, 7).Value = "1"
$ws.Cells.Item(3, 8).Value = "0"
$ws.Cells.Item(3, 9).Value = "83.33"
$ws.Cells.Item(3, 10).Value = "Delhi Capitals"
$ws.Cells.Item(3, 11).Value = "Sharjah"
$ws.Cells.Item(3, 12).Value = "October 02"
$ws.Cells.Item(3, 13).Value = "Capitals won by 4 wickets (with 5 balls remaining)"

# Row 4: 30th, vs Chennai Super Kings
$ws.Cells.Item(4, 1).Value = "30th"
$ws.Cells.Item(4, 2).Value = "Mumbai Indians"
$ws.Cells.Item(4, 3).Value = "Saurabh Tiwary"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = "50"
$ws.Cells.Item(4, 6).Value = "40"
$ws.Cells.Item(4, 7).Value = "5"
$ws.Cells.Item(4, 8).Value = "0"
$ws.Cells.Item(4, 9).Value = "125.00"
$ws.Cells.Item(4, 10).Value = "Chennai Super Kings"
$ws.Cells.Item(4, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(4, 12).Value = "September 19"
$ws.Cells.Item(4, 13).Value = "Super Kings won by 20 runs"

# Row 5: 34th, vs Kolkata Knight Riders
$ws.Cells.Item(5, 1).Value = "34th"
$ws.Cells.Item(5, 2).Value = "Mumbai Indians"
$ws.Cells.Item(5, 3).Value = "Saurabh Tiwary"
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = "5"
$ws.Cells.Item(5, 6).Value = "2"
$ws.Cells.Item(5, 7).Value = "1"
$ws.Cells.Item(5, 8).Value = "0"
$ws.Cells.Item(5, 9).Value = "250.00"
$ws.Cells.Item(5, 10).Value = "Kolkata Knight Riders"
$ws.Cells.Item(5, 11).Value = "Abu Dhabi"
$ws.Cells.Item(5, 12).Value = "September 23"
$ws.Cells.Item(5, 13).Value = "KKR won by 7 wickets (with 29 balls remaining)"
